# Day 4 & 5 added - append 4 new DRS review rows for Match 23 (PBKS vs SRH)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(23, "PBKS", "SRH", 1, "SRH",  "PBKS", 6,  "PBKS", "NA Patwardhan", "NAP", "Wide",   "Called",     "Not Called", "RA Tripathi",    "HV Patel",       "Successful",   "No"),
    @(23, "PBKS", "SRH", 1, "SRH",  "PBKS", 10, "PBKS", "NA Patwardhan", "NAP", "Wicket", "Not Out",    "Out",        "RA Tripathi",    "HV Patel",       "Successful",   "No"),
    @(23, "PBKS", "SRH", 1, "SRH",  "PBKS", 17, "SRH",  "Navdeep Singh", "NS",  "Wide",   "Not Called", "Not Called", "Abdul Samad",    "Arshdeep Singh", "Unsuccessful", "No"),
    @(23, "PBKS", "SRH", 2, "PBKS", "SRH",  17, "SRH",  "Navdeep Singh", "NS",  "Wide",   "Called",     "Called",     "Shashank Singh", "B Kumar",        "Unsuccessful", "No")
)

$startRow = 71
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}
